# add teams names to html web page
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the placeholder team names with the real team names (TeamA / TeamB)
$ws.Range("B2").Value = "Israel"
$ws.Range("C2").Value = "Brazil"

# Update the current selection to match the saved view state
$ws.Range("G10").Select()
